$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: replace "Email" with "Section 1" and add "Section 2"
$ws.Range("E1").Value = "Section 1"
$ws.Range("F1").Value = "Section 2"

# Populate "Section 1" column (E) with 2 for every student row
$ws.Range("E2:E17").Value = 2

# Populate "Section 2" column (F) with 2 for every student row
$ws.Range("F2:F17").Value = 2

# Two students were absent for section 2
$ws.Range("F3").Value = "Abs"
$ws.Range("F11").Value = "Abs"

# Update selection to F3 as in the saved file
$ws.Range("F3").Select()
